# "Made progress on data flow"
#
# Appends one new data row (row 39: B/C/D only, no A) to Sheet1 with the
# values finish_SA_click / aspects / aspects,traits,aspects_traits, using
# the same light-fill "section" style already used by rows 2, 10, 14 and
# 29 (cellXf index 3 in the original workbook). The new strings land at
# the end of the shared-string table, and the sheet's dimension/selection
# are updated to reflect the newly used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row - only columns B:D are populated (A39 stays empty).
$ws.Range("B39").Value = "finish_SA_click"
$ws.Range("C39").Value = "aspects"
$ws.Range("D39").Value = "aspects,traits,aspects_traits"

# Pick up the light-fill "section header" formatting used elsewhere in the
# sheet (e.g. row 10) and apply it to the new row, then drop the leading
# cell (A39) so only B39:D39 end up carrying that style, matching the
# other rows that use it only across their populated cells.
$ws.Range("A10:D10").Copy()
$ws.Range("A39:D39").PasteSpecial(-4122)
$ws.Range("A39").Clear()
$excel.CutCopyMode = $false

# After typing across B39:D39 and confirming, Excel's selection lands on
# the first column of the next row.
$ws.Range("B40").Select() | Out-Null
